# Week 13 logging update
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet: append this week's yardage figures to the running logs
# ---------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = "4 2 4 0 2 7 1 0 3 -1 13 3 0 2 4 1 4 11 1 10 1 1 3 5 -2 5 4 8 16 6 13 5 6 4 5 5 6 19 5 28 3 23 3 1 -1 5 4 3 3 5 9 1 4 4 8 11 7 5 -1 6 5 7 6 4 9 5 8 1 2 10 3 0 0 13 2 -1 5 4 8 3 1 5 4 1 5 -6 2 2 19 0 8 7 5 7 1 4 0 1 2 6 1 11 5 5 47 14 0 14 0 0 10 5 2 18 8 1 5 2 15 4 4 6 6 17 4 5 3 6 6 3 5 12 0 0 14 0 4 3 7 -2 5 3 14 0 10 2 2 7 4 13 2 2 2 3 2 2 5 0 2 -1 9 4 1 12 3 4 -1 0 0 0 1 3 9 0 2 4 4 21 8 8 10 2 7 8 2 6 2 3 5 4 2 7 2 5 2 2 2 6 2 4 2 2 4 2 1 2 1 4 5 8 1 7 11 7 2 3 -2 7 2 5 10 9 3 5 4 -3 4 2 -3 3 2 1 9 9 2 2 2 1 0 3 1 3 15 22 12 7 2 3 -1 2 5 31 3 1 5 1 5 2 8 3 5 4 4 4 4 2 5 8 1 3 5 2 1 4 3 2 3 8 3 1 12 -5 1 1 0 3 2 1 2 33 2 1 3 4 -2 4 58 2 10 2 3 5 3 3 10"
$wsYDS.Range("B3").Value = "28 5 16 11 8 12 7 8 7 22 2 5 22 6 8 3 4 13 -5 4 11 7 10 5 16 6 14 6 4 3 7 10 9 6 21 2 3 6 11 9 31 10 6 7 6 20 7 17 12 10 13 5 16 10 34 4 13 -4 1 11 7 7 5 12 4 8 8 44 7 8 4 12 19 9 8 6 9 6 15 15 13 11 7 5 22 18 12 10 14 18 4 7 20 35 6 7 6 23 8 6 8 6 2 49 6 8 -2 11 21 7 15 24 11 30 15 4 8 14 18 22 35 7 11 14 11 18 17 3 23 1 13 6 13 44 20 11 3 11 4 13 0 18 27 20 1 13 12 33 3 2 4 12 24 14 6 9 7 35 5 8 17 8 20 29 4 9 3 7 73 10 2 35 1 13 15 7 19 33 18 1 7 -3 15 5 5 0 32 17 21 4 18 5 4 24 9 32 5 10 7 4 20 11 4 37 11 13 5 9 12 6 6 18 5 7 4 9 8 16 9 14 9 19 7 22 11 6 8 -3 11 11 2 2 8 11 3 1 5 4 5 16 4 0 6 12 18 14 -2 3 19 5 19 7 11 16 1 6 15 17 10 10 9 4 51 5 13 5 7 7 15 7 41 12 7 6 5 4 32 32 11 5 3 17 10 2 2 4 9 5 1 41 7 1 21 10 10 5 17 1 4 20 0 12 2 10 1 3 4 13 25 9 4 3"
$wsYDS.Range("C2").Value = "5 3 5 3 3 5 2 5 1 7 0 6 6 1 20 -1 16 1 3 1 5 1 0 5 1 9 5 3 9 3 7 4 -2 24 2 3 3 1 3 7 9 8 2 4 14 12 2 -1 1 -1 14 7 4 11 -2 2 6 7 3 3 0 7 -2 7 7 4 6 5 1 4 3 4 1 14 2 4 1 0 0 -1 4 0 2 4 -2 7 3 4 21 4 4 5 3 1 1 1 11 -1 6 9 1 3 4 1 21 6 6 3 1 0 0 6 5 2 6 2 -1 15 0 2 3 2 0 0 13 19 4 3 -2 8 4 8 8 7 1 4 -4 -2 11 20 -3 7 13 1 -2 3 17 0 12 2 2 11 2 2 2 30 2 1 1 7 14 6 7 6 4 -4 2 0 1 10 3 2 0 2 1 -1 6 14 3 2 3 5 -3 4 1 -5 1 6 2 4 7 17 6 8 33 11 4 3 4 1 5 1 5 0 2 3 3 7 4 14 3 6 4 10 -2 -7 3 2 2 4 1 3 0 9 4 1 2 3 0 10 5 5 4 -4 6 6 3 22 2 4 4 4 2 4 5 21 -1 0 4 4 2 5 3 -1 3 2 9 9 3 -1 25 24 7 0 4 8 6 5 9 23 4 6 8 0 7 7"
$wsYDS.Range("C3").Value = "10 19 28 16 8 5 27 5 19 17 2 10 47 20 2 2 11 10 1 5 7 9 7 9 7 3 11 8 4 6 20 24 22 11 7 0 5 14 -2 12 9 27 12 14 8 -1 6 10 7 8 42 8 -2 5 7 16 8 20 4 2 15 30 14 38 5 5 9 19 7 19 28 41 27 0 3 4 1 5 12 5 7 24 19 16 12 15 -4 9 9 29 9 3 13 6 5 17 1 8 5 9 39 9 15 2 11 4 7 13 8 6 5 5 8 38 5 9 18 28 6 35 6 2 26 5 13 13 16 14 13 13 11 7 13 9 14 22 20 27 0 8 7 3 5 11 3 11 75 14 31 3 20 3 5 7 -2 6 11 18 7 7 -4 -2 4 4 9 -3 6 13 6 3 -3 25 5 8 11 44 19 7 6 9 17 6 3 11 20 8 -3 16 40 17 18 17 14 0 8 13 21 9 1 10 8 12 0 8 37 13 23 7 11 8 24 14 22 10 12 6 6 4 8 8 8 12 56 5 12 12 21 17 13 31 3 11 -1 4 35 10 -5 16 54 9 9 4 30 7 8 6 12 -2 -6 25 27 9 24 9 1 9 30 6 1 5 13 13 70"

# ---------------------------------------------------------------
# OFF sheet: Home (row 2) / Road (row 3) season totals
# ---------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 159
$wsOFF.Range("F2").Value = 43
$wsOFF.Range("G2").Value = 40
$wsOFF.Range("N2").Value = 12
$wsOFF.Range("O2").Value = 18

$wsOFF.Range("B3").Value = 12
$wsOFF.Range("C3").Value = 131
$wsOFF.Range("D3").Value = 6
$wsOFF.Range("E3").Value = 28
$wsOFF.Range("F3").Value = 96
$wsOFF.Range("G3").Value = 31
$wsOFF.Range("I3").Value = 45
$wsOFF.Range("J3").Value = 44
$wsOFF.Range("L3").Value = 262
$wsOFF.Range("M3").Value = 182
$wsOFF.Range("Q3").Value = 420

# ---------------------------------------------------------------
# DEF sheet: Home (row 2) / Road (row 3) season totals
# ---------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 137
$wsDEF.Range("E2").Value = 11
$wsDEF.Range("F2").Value = 41
$wsDEF.Range("G2").Value = 34
$wsDEF.Range("I2").Value = 5
$wsDEF.Range("J2").Value = 17
$wsDEF.Range("N2").Value = 10
$wsDEF.Range("O2").Value = 13
$wsDEF.Range("P2").Value = 8

$wsDEF.Range("B3").Value = 15
$wsDEF.Range("C3").Value = 144
$wsDEF.Range("F3").Value = 80
$wsDEF.Range("G3").Value = 27
$wsDEF.Range("H3").Value = 21
$wsDEF.Range("I3").Value = 42
$wsDEF.Range("J3").Value = 41
$wsDEF.Range("L3").Value = 225
$wsDEF.Range("M3").Value = 143
$wsDEF.Range("Q3").Value = 370

# ---------------------------------------------------------------
# ST sheet: special-teams totals + per-kicker/punter logs
# ---------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value = 72
$wsST.Range("D2").Value = 44
$wsST.Range("F2").Value = 339
$wsST.Range("G2").Value = 327
$wsST.Range("J2").Value = 164
$wsST.Range("K2").Value = 155
$wsST.Range("N2").Value = 68
$wsST.Range("O2").Value = 38

$wsST.Range("B3").Value = 61

$wsST.Range("D3").Value = "38 47 62 44 56 49 40 57 51 40 36 63 61 51 50 54 36 41 55 55 46 38 51 60 52 46 47 34 48 53 60 55 53 44 58 56 59 60 59 59 54 31 37 37"
$wsST.Range("D4").Value = "0 24 7 6 -4 -3 0 5 0 0 0 0 7 23 0 7 0 0 0 0 1 0 5 4 0 0 0 0 0 9 21 0 0 0 0 9 7 8 19 0 8 0 2 0"
$wsST.Range("D5").Value = "0 3 12 0 5 8 0 0 0 0 0 0 0 0 17 0 -1 0 -1 0 0 7 0 -3 0 0 0 3 16 0 0 7 0 0 2 0 0 0 14 7 16 9 0 0 12"
$wsST.Range("B6").Value = "24 23 24 14 25 22 27 22 54 35 29 21 18 18 0 11 22"

# ---------------------------------------------------------------
# TURNS sheet: Road (row 3) turnover totals
# ---------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("B3").Value = 7
$wsTURNS.Range("C3").Value = 9
$wsTURNS.Range("E3").Value = 4

# ---------------------------------------------------------------
# PEN sheet: Holding (row 3) count
# ---------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B3").Value = 24
